# Updated cryptos list with GitHub Actions: refresh Price (D) and
# Volume(1h) (E) columns on the crypto table in Sheet1.
#
# Some Price values (column D) are plain decimal numbers (e.g. "214.51").
# The source data stores these as literal text (inlineStr) rather than
# numeric cells, so we force the cell's number format to Text ("@") before
# assigning the value - this keeps Excel from auto-converting the literal
# into a floating point number (which would corrupt values like "15.90"
# by dropping the trailing zero, or introduce binary float noise).
# Values that contain two dots (e.g. "26.907.58") are never parsed as
# numbers by Excel, so no such guard is needed for those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.907.58"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "1.668.86"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.51"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.25"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").Value = "1.904.35"
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "1.662.95"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.45"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "26.917.86"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "233.16"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("E20").Value = "  +0.25%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.14"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("E24").Value = "  -2.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.33"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.90"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -2.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0497"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.32"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").Value = "1.454.50"
$ws.Range("E33").Value = "  -6.20%  "
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("E35").Value = "  +2.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.583"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("E40").Value = "  +14.45%  "
$ws.Range("E41").Value = "  -4.48%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.18"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").Value = "1.811.92"
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.779"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.61"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.59"
$ws.Range("E51").Value = "  -0.52%  "
